$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a second row of test-case data (row 3), mirroring the existing
# "sample1 / NavigateTo / NA / google.com" row above it.
# Order matches the shared-string insertion order of the target workbook:
# sample2, "2" (text/quote-prefixed SNO), the bing URL, then the keyword.
$ws.Range("B3").Value = "sample2"
$ws.Range("A3").Value = "'2"
$ws.Range("E3").Value = "https://www.bing.com/"
$ws.Range("C3").Value = "explicitwait"
$ws.Range("D3").Value = "NA"

# Turn the TEST DATA cell into a real hyperlink, just like E2, then restore
# the shared "Hyperlink" cell style so it renders identically to row 2.
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.bing.com/") | Out-Null
$ws.Range("E3").Style = "Hyperlink"

# Recorded UI selection state after the edit.
$ws.Range("D6").Select() | Out-Null
